$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q3" right after "总计" (sheet 1),
#    pushing every other quarterly sheet down by one position.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $total)
$q3.Name = "2022-Q3"

# Header row (copy style from the "总计" sheet header, which already has the
# bordered/bold/centered style used for every sheet's header row).
$total.Cells.Item(1, 2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, $c + 2).Value2 = $headers[$c]
}

# Force the fund-code / fund-name / size / position / ratio / value columns to be
# stored as text (matches the source data, which keeps things like leading zeros
# in fund codes and fixed-decimal strings such as "29.98").
$q3.Range("B2:G38").NumberFormat = "@"

$q3rows = @(
    ,@("010488","鹏华优选成长混合A","29.98","67.49","3.37","1.0103",6)
    ,@("001230","鹏华医药科技股票","15.31","80.22","5.32","0.8145",5)
    ,@("011568","鹏华产业升级混合A","18.20","76.63","3.54","0.6443",7)
    ,@("000242","景顺长城策略精选","11.62","90.68","4.63","0.5380",6)
    ,@("501011","汇添富中证中药指数（LOF）A","11.36","94.73","3.61","0.4101",7)
    ,@("003291","信澳健康中国灵活配置混合A","7.25","91.02","5.24","0.3799",3)
    ,@("160610","鹏华动力增长混合（LOF）","13.14","53.95","2.55","0.3351",6)
    ,@("000780","鹏华医疗保健股票","6.60","81.50","4.63","0.3056",6)
    ,@("005505","前海开源中药研究精选股票A","4.10","83.96","7.40","0.3034",6)
    ,@("011598","信澳医药健康混合","5.45","88.25","5.28","0.2878",3)
    ,@("012093","鹏华创新升级混合A","5.82","64.25","4.80","0.2794",3)
    ,@("005506","前海开源中药研究精选股票C","3.58","83.96","7.40","0.2649",6)
    ,@("501012","汇添富中证中药指数（LOF）C","6.42","94.73","3.61","0.2318",7)
    ,@("159647","鹏华中证中药ETF","6.16","94.79","3.68","0.2267",6)
    ,@("015208","信澳健康中国灵活配置混合C","3.72","91.02","5.24","0.1949",3)
    ,@("003713","英大睿盛灵活配置混合A","2.83","93.65","6.88","0.1947",6)
    ,@("003714","英大睿盛灵活配置混合C","2.19","93.65","6.88","0.1507",6)
    ,@("160603","鹏华普天收益混合","4.25","65.73","2.97","0.1262",6)
    ,@("001678","英大国企改革主题股票","1.55","93.30","7.35","0.1139",1)
    ,@("562390","银华中证中药ETF","2.34","98.09","3.73","0.0873",7)
    ,@("561510","华泰柏瑞中证中药ETF","2.02","95.98","3.72","0.0751",6)
    ,@("002259","鹏华健康环保灵活配置混合","1.96","77.62","3.41","0.0668",6)
    ,@("510081","长盛动态精选混合","2.73","57.03","2.42","0.0661",7)
    ,@("011331","鹏华远见成长混合A","1.81","68.90","3.31","0.0599",6)
    ,@("001524","华泰柏瑞精选回报灵活配置混合","5.56","21.47","1.06","0.0589",6)
    ,@("001607","英大策略优选混合A","0.57","91.98","6.86","0.0391",3)
    ,@("010489","鹏华优选成长混合C","0.79","67.49","3.37","0.0266",6)
    ,@("012522","英大稳固增强核心一年持有混合C","1.24","27.71","1.53","0.0190",7)
    ,@("003447","英大睿鑫灵活配置混合C","0.21","92.71","7.85","0.0165",4)
    ,@("012521","英大稳固增强核心一年持有混合A","0.75","27.71","1.53","0.0115",7)
    ,@("011569","鹏华产业升级混合C","0.28","76.63","3.54","0.0099",7)
    ,@("012094","鹏华创新升级混合C","0.17","64.25","4.80","0.0082",3)
    ,@("011179","浙商智选食品饮料股票A","0.14","91.42","5.80","0.0081",9)
    ,@("011332","鹏华远见成长混合C","0.19","68.90","3.31","0.0063",6)
    ,@("003446","英大睿鑫灵活配置混合A","0.07","92.71","7.85","0.0055",4)
    ,@("011180","浙商智选食品饮料股票C","0.08","91.42","5.80","0.0046",9)
    ,@("001608","英大策略优选混合C","0.02","91.98","6.86","0.0014",3)
)

$r = 2
foreach ($row in $q3rows) {
    $q3.Cells.Item($r, 1).Value2 = $r - 2
    $q3.Cells.Item($r, 2).Value2 = $row[0]
    $q3.Cells.Item($r, 3).Value2 = $row[1]
    $q3.Cells.Item($r, 4).Value2 = $row[2]
    $q3.Cells.Item($r, 5).Value2 = $row[3]
    $q3.Cells.Item($r, 6).Value2 = $row[4]
    $q3.Cells.Item($r, 7).Value2 = $row[5]
    $q3.Cells.Item($r, 8).Value2 = $row[6]
    $r++
}

# Column A (row index, 0-based) keeps the bordered/centered numeric style used by
# every other sheet's "A" column.
$total.Cells.Item(2, 1).Copy()
$q3.Range("A2:A38").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new top data row for 2022-Q3
#    and push the existing quarters down by one row.
# ---------------------------------------------------------------------------
for ($row = 8; $row -ge 2; $row--) {
    $dst = $row + 1
    $total.Cells.Item($dst, 2).Value2 = $total.Cells.Item($row, 2).Value2
    $total.Cells.Item($dst, 3).Value2 = $total.Cells.Item($row, 3).Value2
    $total.Cells.Item($dst, 4).Value2 = $total.Cells.Item($row, 4).Value2
}

# Row 9 is brand new ("2020-Q4" used to live in row 8) - copy the numeric-index
# style from row 8's A cell before writing the shifted index value into it.
$total.Cells.Item(8, 1).Copy()
$total.Cells.Item(9, 1).PasteSpecial(-4122)
$total.Cells.Item(9, 1).Value2 = 7

$total.Cells.Item(2, 2).Value2 = "2022-Q3"
$total.Cells.Item(2, 3).Value2 = 37
$total.Cells.Item(2, 4).Value2 = 7.38
